{"js": "// Update the worksheet date and the 25 two-digit-by-two-digit multiplication\n// problems in place. Each old value appears exactly once in the document, so\n// a simple search + insertText(\"Replace\") pass for every pair is both safe\n// and deterministic.\n\nconst pairs = [\n  [\"2024-04-13 Saturday\", \"2024-04-14 Sunday\"],\n  [\"70\u00d789=\", \"66\u00d792=\"],\n  [\"29\u00d740=\", \"36\u00d751=\"],\n  [\"97\u00d765=\", \"18\u00d740=\"],\n  [\"84\u00d713=\", \"91\u00d731=\"],\n  [\"98\u00d751=\", \"53\u00d716=\"],\n  [\"78\u00d749=\", \"11\u00d722=\"],\n  [\"69\u00d721=\", \"49\u00d727=\"],\n  [\"75\u00d762=\", \"72\u00d786=\"],\n  [\"90\u00d780=\", \"59\u00d757=\"],\n  [\"57\u00d729=\", \"71\u00d740=\"],\n  [\"45\u00d770=\", \"20\u00d726=\"],\n  [\"35\u00d766=\", \"50\u00d756=\"],\n  [\"32\u00d759=\", \"54\u00d771=\"],\n  [\"93\u00d732=\", \"91\u00d793=\"],\n  [\"25\u00d717=\", \"96\u00d789=\"],\n  [\"34\u00d760=\", \"44\u00d787=\"],\n  [\"70\u00d794=\", \"38\u00d732=\"],\n  [\"12\u00d745=\", \"74\u00d789=\"],\n  [\"93\u00d795=\", \"55\u00d773=\"],\n  [\"47\u00d745=\", \"88\u00d725=\"],\n  [\"57\u00d776=\", \"20\u00d790=\"],\n  [\"61\u00d785=\", \"15\u00d752=\"],\n  [\"29\u00d713=\", \"37\u00d766=\"],\n  [\"14\u00d740=\", \"74\u00d731=\"],\n  [\"95\u00d753=\", \"97\u00d799=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of pairs) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  // Each value is unique in this document, so exactly one hit is expected;\n  // replace only the first occurrence found to stay a precise 1:1 swap.\n  if (results.items.length > 0) {\n    results.items[0].insertText(newText, \"Replace\");\n    await context.sync();\n  }\n}\n", "ps1": "# Update the worksheet date and the 25 two-digit-by-two-digit multiplication\n# problems in place. Each old value appears exactly once in the document, so\n# a simple ordered Find/Replace pass (wdReplaceOne) for every pair is both\n# safe and deterministic.\n\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"2024-04-13 Saturday\", \"2024-04-14 Sunday\"),\n    @(\"70\u00d789=\", \"66\u00d792=\"),\n    @(\"29\u00d740=\", \"36\u00d751=\"),\n    @(\"97\u00d765=\", \"18\u00d740=\"),\n    @(\"84\u00d713=\", \"91\u00d731=\"),\n    @(\"98\u00d751=\", \"53\u00d716=\"),\n    @(\"78\u00d749=\", \"11\u00d722=\"),\n    @(\"69\u00d721=\", \"49\u00d727=\"),\n    @(\"75\u00d762=\", \"72\u00d786=\"),\n    @(\"90\u00d780=\", \"59\u00d757=\"),\n    @(\"57\u00d729=\", \"71\u00d740=\"),\n    @(\"45\u00d770=\", \"20\u00d726=\"),\n    @(\"35\u00d766=\", \"50\u00d756=\"),\n    @(\"32\u00d759=\", \"54\u00d771=\"),\n    @(\"93\u00d732=\", \"91\u00d793=\"),\n    @(\"25\u00d717=\", \"96\u00d789=\"),\n    @(\"34\u00d760=\", \"44\u00d787=\"),\n    @(\"70\u00d794=\", \"38\u00d732=\"),\n    @(\"12\u00d745=\", \"74\u00d789=\"),\n    @(\"93\u00d795=\", \"55\u00d773=\"),\n    @(\"47\u00d745=\", \"88\u00d725=\"),\n    @(\"57\u00d776=\", \"20\u00d790=\"),\n    @(\"61\u00d785=\", \"15\u00d752=\"),\n    @(\"29\u00d713=\", \"37\u00d766=\"),\n    @(\"14\u00d740=\", \"74\u00d731=\"),\n    @(\"95\u00d753=\", \"97\u00d799=\")\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 1)\n}\n"}
